$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H107").Value = 333.9
$ws_ALC.Range("I107").Value = 287.33334
$ws_ALC.Range("K107").Value = 287.33334
$ws_ALC.Range("M107").Value = 1632.66666

$ws_ALC.Range("H129").Value = 1036.9333
$ws_ALC.Range("J129").Value = 1127.3334
$ws_ALC.Range("L129").Value = 3382.0002
$ws_ALC.Range("N129").Value = -13382.0002

$ws_ALC.Range("H132").Value = 3292870.5
$ws_ALC.Range("I132").Value = 3427892.5
$ws_ALC.Range("J132").Value = 7335.3335
$ws_ALC.Range("K132").Value = 10283677.5
$ws_ALC.Range("L132").Value = 22006.0005
$ws_ALC.Range("M132").Value = -10281147.5
$ws_ALC.Range("N132").Value = -27066.0005

$ws_ALC.Range("H137").Value = 1598.2554
$ws_ALC.Range("I137").Value = 1097.9269
$ws_ALC.Range("J137").Value = 5017.1665
$ws_ALC.Range("K137").Value = 3293.7807
$ws_ALC.Range("L137").Value = 15051.4995
$ws_ALC.Range("M137").Value = -743.7806999999998
$ws_ALC.Range("N137").Value = -20151.4995

$ws_ALC.Range("H138").Value = 5360.974
$ws_ALC.Range("I138").Value = 1257.8857
$ws_ALC.Range("J138").Value = 8780.214
$ws_ALC.Range("K138").Value = 3773.6571
$ws_ALC.Range("L138").Value = 26340.642
$ws_ALC.Range("M138").Value = 1366.3429
$ws_ALC.Range("N138").Value = -36620.642

$ws_ALC.Range("H141").Value = 1593.1333
$ws_ALC.Range("I141").Value = 1360.5
$ws_ALC.Range("J141").Value = 4850
$ws_ALC.Range("K141").Value = 4081.5
$ws_ALC.Range("L141").Value = 14550
$ws_ALC.Range("M141").Value = 1098.5
$ws_ALC.Range("N141").Value = -24910

$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H2").Value = 72367.5
$ws_ARM.Range("I2").Value = 925.7778
$ws_ARM.Range("J2").Value = 200962.6
$ws_ARM.Range("K2").Value = 925.7778
$ws_ARM.Range("L2").Value = 200962.6
$ws_ARM.Range("M2").Value = -812.7778
$ws_ARM.Range("N2").Value = -201188.6

$ws_ARM.Range("H6").Value = 14383.917
$ws_ARM.Range("I6").Value = 27501
$ws_ARM.Range("J6").Value = 11760.5
$ws_ARM.Range("K6").Value = 27501
$ws_ARM.Range("L6").Value = 11760.5
$ws_ARM.Range("M6").Value = -27328
$ws_ARM.Range("N6").Value = -12106.5

$ws_ARM.Range("H16").Value = 5333.3335
$ws_ARM.Range("I16").Value = 4000
$ws_ARM.Range("J16").Value = 8000
$ws_ARM.Range("K16").Value = 4000
$ws_ARM.Range("L16").Value = 8000
$ws_ARM.Range("M16").Value = -3713
$ws_ARM.Range("N16").Value = -8574

$ws_ARM.Range("H32").Value = 22255.592
$ws_ARM.Range("I32").Value = 3529.581
$ws_ARM.Range("J32").Value = 220216.28
$ws_ARM.Range("K32").Value = 3529.581
$ws_ARM.Range("L32").Value = 220216.28
$ws_ARM.Range("M32").Value = -3242.581
$ws_ARM.Range("N32").Value = -220790.28

$ws_ARM.Range("H110").Value = 58942350
$ws_ARM.Range("I110").Value = 83500980
$ws_ARM.Range("J110").Value = 1652.6
$ws_ARM.Range("K110").Value = 83500980
$ws_ARM.Range("L110").Value = 1652.6
$ws_ARM.Range("M110").Value = -83498935
$ws_ARM.Range("N110").Value = -5742.6

$ws_ARM.Range("H116").Value = 72367.5
$ws_ARM.Range("I116").Value = 925.7778
$ws_ARM.Range("J116").Value = 200962.6
$ws_ARM.Range("K116").Value = 925.7778
$ws_ARM.Range("L116").Value = 200962.6
$ws_ARM.Range("M116").Value = 1368.2222
$ws_ARM.Range("N116").Value = -205550.6

$ws_ARM.Range("H132").Value = 1680.64
$ws_ARM.Range("I132").Value = 1500.7246
$ws_ARM.Range("J132").Value = 3749.6667
$ws_ARM.Range("K132").Value = 4502.1738
$ws_ARM.Range("L132").Value = 11249.0001
$ws_ARM.Range("M132").Value = -1972.1738
$ws_ARM.Range("N132").Value = -16309.0001

$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H3").Value = 72367.5
$ws_BSM.Range("I3").Value = 925.7778
$ws_BSM.Range("J3").Value = 200962.6
$ws_BSM.Range("K3").Value = 925.7778
$ws_BSM.Range("L3").Value = 200962.6
$ws_BSM.Range("M3").Value = -811.7778
$ws_BSM.Range("N3").Value = -201190.6

$ws_BSM.Range("H94").Value = 699.7646999999999
$ws_BSM.Range("I94").Value = 695.2
$ws_BSM.Range("J94").Value = 701.6667
$ws_BSM.Range("K94").Value = 695.2
$ws_BSM.Range("L94").Value = 701.6667
$ws_BSM.Range("M94").Value = -244.2
$ws_BSM.Range("N94").Value = -1603.6667

$ws_BSM.Range("H107").Value = 100086104
$ws_BSM.Range("I107").Value = 125107464
$ws_BSM.Range("J107").Value = 657
$ws_BSM.Range("K107").Value = 125107464
$ws_BSM.Range("L107").Value = 657
$ws_BSM.Range("M107").Value = -125105544
$ws_BSM.Range("N107").Value = -4497

$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H16").Value = 1319.75
$ws_CRP.Range("I16").Value = 1223.75
$ws_CRP.Range("J16").Value = 1415.75
$ws_CRP.Range("K16").Value = 1223.75
$ws_CRP.Range("L16").Value = 1415.75
$ws_CRP.Range("M16").Value = -936.75
$ws_CRP.Range("N16").Value = -1989.75

$ws_CRP.Range("H68").Value = 17571.773
$ws_CRP.Range("J68").Value = 17571.773
$ws_CRP.Range("L68").Value = 17571.773
$ws_CRP.Range("N68").Value = -19069.773

$ws_CRP.Range("H71").Value = 17571.773
$ws_CRP.Range("J71").Value = 17571.773
$ws_CRP.Range("L71").Value = 52715.319
$ws_CRP.Range("N71").Value = -60203.319

$ws_CRP.Range("H107").Value = 9096.416999999999
$ws_CRP.Range("J107").Value = 619.5
$ws_CRP.Range("L107").Value = 619.5
$ws_CRP.Range("N107").Value = -4459.5

$ws_CRP.Range("H113").Value = 1319.75
$ws_CRP.Range("I113").Value = 1223.75
$ws_CRP.Range("J113").Value = 1415.75
$ws_CRP.Range("K113").Value = 1223.75
$ws_CRP.Range("L113").Value = 1415.75
$ws_CRP.Range("M113").Value = 946.25
$ws_CRP.Range("N113").Value = -5755.75

$ws_CRP.Range("H132").Value = 23079184
$ws_CRP.Range("I132").Value = 20002180
$ws_CRP.Range("K132").Value = 60006540
$ws_CRP.Range("M132").Value = -60004010

$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H4").Value = 90909220
$ws_CUL.Range("I4").Value = 139.9
$ws_CUL.Range("J4").Value = 1000000000
$ws_CUL.Range("K4").Value = 419.7
$ws_CUL.Range("L4").Value = 3000000000
$ws_CUL.Range("M4").Value = -307.7
$ws_CUL.Range("N4").Value = -3000000224

$ws_CUL.Range("H5").Value = 1640.3636
$ws_CUL.Range("I5").Value = 836.9286
$ws_CUL.Range("J5").Value = 2232.3684
$ws_CUL.Range("K5").Value = 2510.7858
$ws_CUL.Range("L5").Value = 6697.1052
$ws_CUL.Range("M5").Value = -2398.7858
$ws_CUL.Range("N5").Value = -6921.1052

$ws_CUL.Range("H114").Value = 753.8823
$ws_CUL.Range("I114").Value = 280.3
$ws_CUL.Range("J114").Value = 1430.4286
$ws_CUL.Range("K114").Value = 840.9000000000001
$ws_CUL.Range("L114").Value = 4291.2858
$ws_CUL.Range("M114").Value = 2413.1
$ws_CUL.Range("N114").Value = -10799.2858

$ws_CUL.Range("H121").Value = 7787.65
$ws_CUL.Range("I121").Value = 5429.8887
$ws_CUL.Range("J121").Value = 9716.727999999999
$ws_CUL.Range("K121").Value = 16289.6661
$ws_CUL.Range("L121").Value = 29150.184
$ws_CUL.Range("M121").Value = -14979.6661
$ws_CUL.Range("N121").Value = -31770.184

$ws_CUL.Range("H131").Value = 870.35
$ws_CUL.Range("I131").Value = 725
$ws_CUL.Range("J131").Value = 873.3163500000001
$ws_CUL.Range("K131").Value = 2175
$ws_CUL.Range("L131").Value = 2619.94905
$ws_CUL.Range("M131").Value = 2865
$ws_CUL.Range("N131").Value = -12699.94905

$ws_CUL.Range("H135").Value = 1640.3636
$ws_CUL.Range("I135").Value = 836.9286
$ws_CUL.Range("J135").Value = 2232.3684
$ws_CUL.Range("K135").Value = 7532.3574
$ws_CUL.Range("L135").Value = 20091.3156
$ws_CUL.Range("M135").Value = -4997.3574
$ws_CUL.Range("N135").Value = -25161.3156

$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H5").Value = 2503000
$ws_GSM.Range("J5").Value = 6000
$ws_GSM.Range("L5").Value = 6000
$ws_GSM.Range("N5").Value = -6224

$ws_GSM.Range("H12").Value = 8126125.5
$ws_GSM.Range("I12").Value = 9286286
$ws_GSM.Range("J12").Value = 5000
$ws_GSM.Range("K12").Value = 9286286
$ws_GSM.Range("L12").Value = 5000
$ws_GSM.Range("M12").Value = -9286146
$ws_GSM.Range("N12").Value = -5280

$ws_GSM.Range("H97").Value = 50001836
$ws_GSM.Range("I97").Value = 71430640
$ws_GSM.Range("J97").Value = 1301.6666
$ws_GSM.Range("K97").Value = 71430640
$ws_GSM.Range("L97").Value = 1301.6666
$ws_GSM.Range("M97").Value = -71430144
$ws_GSM.Range("N97").Value = -2293.6666

$ws_GSM.Range("H107").Value = 674213.25
$ws_GSM.Range("I107").Value = 544.2222
$ws_GSM.Range("K107").Value = 544.2222
$ws_GSM.Range("M107").Value = 1375.7778

$ws_GSM.Range("H113").Value = 2300
$ws_GSM.Range("I113").Value = 2300
$ws_GSM.Range("J113").Value = 0
$ws_GSM.Range("K113").Value = 2300
$ws_GSM.Range("L113").Value = 0
$ws_GSM.Range("M113").ClearContents()
$ws_GSM.Range("N113").Value = -130

$ws_GSM.Range("H122").Value = 2202.3333
$ws_GSM.Range("I122").Value = 3701.75
$ws_GSM.Range("J122").Value = 1452.625
$ws_GSM.Range("K122").Value = 11105.25
$ws_GSM.Range("L122").Value = 4357.875
$ws_GSM.Range("M122").Value = -8655.25
$ws_GSM.Range("N122").Value = -9257.875

$ws_GSM.Range("H132").Value = 2130.6487
$ws_GSM.Range("I132").Value = 1297.8334
$ws_GSM.Range("J132").Value = 3668.1538
$ws_GSM.Range("K132").Value = 3893.5002
$ws_GSM.Range("L132").Value = 11004.4614
$ws_GSM.Range("M132").Value = -1363.5002
$ws_GSM.Range("N132").Value = -16064.4614

$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H2").Value = 500000
$ws_LTW.Range("J2").Value = 0
$ws_LTW.Range("L2").Value = 0
$ws_LTW.Range("N2").ClearContents()

$ws_LTW.Range("H61").Value = 1888.3334
$ws_LTW.Range("I61").Value = 870
$ws_LTW.Range("K61").Value = 870
$ws_LTW.Range("M61").Value = -668

$ws_LTW.Range("H82").Value = 1339.75
$ws_LTW.Range("I82").Value = 970.7143
$ws_LTW.Range("J82").Value = 1538.4615
$ws_LTW.Range("K82").Value = 970.7143
$ws_LTW.Range("L82").Value = 1538.4615
$ws_LTW.Range("M82").Value = -609.7143
$ws_LTW.Range("N82").Value = -2260.4615

$ws_LTW.Range("H85").Value = 1339.75
$ws_LTW.Range("I85").Value = 970.7143
$ws_LTW.Range("J85").Value = 1538.4615
$ws_LTW.Range("K85").Value = 970.7143
$ws_LTW.Range("L85").Value = 1538.4615
$ws_LTW.Range("M85").Value = 277.2857
$ws_LTW.Range("N85").Value = -4034.4615

$ws_LTW.Range("H113").Value = 1888.3334
$ws_LTW.Range("I113").Value = 870
$ws_LTW.Range("K113").Value = 870
$ws_LTW.Range("M113").Value = 1300

$ws_LTW.Range("H132").Value = 2127.3774
$ws_LTW.Range("I132").Value = 2245.7708
$ws_LTW.Range("J132").Value = 990.8
$ws_LTW.Range("K132").Value = 6737.312399999999
$ws_LTW.Range("L132").Value = 2972.4
$ws_LTW.Range("M132").Value = -4207.312399999999
$ws_LTW.Range("N132").Value = -8032.4

$ws_LTW.Range("H136").Value = 1217.4193
$ws_LTW.Range("I136").Value = 1217.4193
$ws_LTW.Range("J136").Value = 0
$ws_LTW.Range("K136").Value = 3652.2579
$ws_LTW.Range("L136").Value = 0
$ws_LTW.Range("M136").ClearContents()
$ws_LTW.Range("N136").Value = -1102.2579

$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H2").Value = 779753.4399999999
$ws_WVR.Range("I2").Value = 840582.9399999999
$ws_WVR.Range("J2").Value = 49800
$ws_WVR.Range("K2").Value = 840582.9399999999
$ws_WVR.Range("L2").Value = 49800
$ws_WVR.Range("M2").Value = -840470.9399999999
$ws_WVR.Range("N2").Value = -50024

$ws_WVR.Range("H132").Value = 3614.36
$ws_WVR.Range("I132").Value = 4310.647
$ws_WVR.Range("J132").Value = 2134.75
$ws_WVR.Range("K132").Value = 12931.941
$ws_WVR.Range("L132").Value = 6404.25
$ws_WVR.Range("M132").Value = -10401.941
$ws_WVR.Range("N132").Value = -11464.25
